$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -74.8554
$ws.Range("B2").Value = -74.7691

$ws.Range("A3").Value = 39.2359
$ws.Range("B3").Value = 39.3019

$ws.Range("A4").Value = -74.0187
$ws.Range("B4").Value = -74.1059

$ws.Range("A5").Value = 39.8716
$ws.Range("B5").Value = 39.8058
